# Populate Point Group (B) and Energy (C) columns for fullerene indices 9-200
# (rows 10-201) on the "C50" worksheet. These rows previously only had the
# "#" index filled in column A; the rest of the per-isomer data (point
# group + energy, in eV) was missing and is being added here.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("C50")

$data = @"
10	C1	28.216
11	C1	29.848
12	C1	29.774
13	C1	28.626
14	C2v	32.622
15	C1	28.986
16	C1	29.092
17	C1	29.659
18	C1	29.101
19	C2	30.016
20	C1	27.725
21	C1	27.182
22	C1	27.308
23	C1	28.08
24	C1	26.886
25	C2	28.144
26	C1	27.939
27	C1	27.968
28	C2	27.259
29	C1	28.574
30	C1	26.524
31	C1	26.677
32	C1	26.274
33	Cs	27.005
34	Cs	26.837
35	C1	26.893
36	C1	26.741
37	C1	27.059
38	C1	26.383
39	C1	27.265
40	C1	25.623
41	C1	25.557
42	C1	26.809
43	Cs	25.409
44	C2v	26.414
45	C1	24.254
46	C1	27.171
47	C1	26.707
48	C1	25.71
49	C1	26.625
50	C1	25.523
51	C1	27.012
52	C1	25.541
53	C1	25.149
54	C1	26.017
55	C1	25.537
56	C1	25.252
57	C1	29.391
58	C1	25.534
59	C1	25.753
60	C1	30.24
61	C1	28.401
62	C2	28.361
63	C1	25.875
64	C1	26.806
65	C2	27.497
66	C1	26.227
67	C1	26.594
68	C2	27.783
69	C1	25.8
70	Cs	25.369
71	C1	25.501
72	C1	27.783
73	Cs	26.034
74	C1	26.938
75	C2	26.73
76	C1	26.376
77	C1	26.949
78	C1	26.327
79	C2	27.894
80	C1	27.873
81	C2	27.815
82	C1	26.611
83	C1	25.64
84	C1	25.279
85	C1	26.596
86	C1	25.717
87	Cs	27.351
88	C1	26.083
89	Cs	25.345
90	C1	26.873
91	C1	26.257
92	Cs	25.234
93	Cs	26.354
94	C2	27.309
95	C1	26.629
96	C2	25.919
97	C1	25.813
98	C1	25.176
99	C1	24.893
100	C1	24.992
101	Cs	24.873
102	C1	25.125
103	C1	24.462
104	C1	25.7
105	C1	24.646
106	C1	25.118
107	C1	24.968
108	C1	24.179
109	C2	25.171
110	C1	24.94
111	C1	24.752
112	C1	24.314
113	C1	25.247
114	C1	24.973
115	C1	24.398
116	C1	24.085
117	C2	26.571
118	C1	24.488
119	C1	24.85
120	C1	24.514
121	C1	25.293
122	C1	25.05
123	C1	25.583
124	C1	25.107
125	C1	25.419
126	C2v	28.328
127	C1	26.433
128	C1	25.388
129	C1	26.236
130	C1	25.401
131	C1	25.496
132	C1	25.835
133	C1	24.932
134	C1	25.397
135	C1	26.498
136	C1	26.311
137	C1	26.851
138	C1	25.558
139	C2	28
140	C1	27.209
141	C1	26.674
142	C1	27.933
143	C1	25.914
144	C1	24.222
145	C1	24.157
146	C1	24.588
147	C1	24.185
148	C1	25.322
149	C1	27.548
150	C1	25.128
151	C1	24.175
152	C1	24.418
153	C1	25.643
154	C1	25.038
155	C1	25.125
156	C1	26.479
157	C2	25.427
158	C3v	24.445
159	Cs	24.481
160	C1	25.427
161	C1	24.879
162	C1	26.202
163	Cs	26.156
164	C1	24.521
165	C1	26.204
166	C1	24.755
167	C1	25.045
168	C2	24.941
169	Cs	24.603
170	C1	23.885
171	C1	25.732
172	C1	23.93
173	C1	24.099
174	C1	23.859
175	C1	24.032
176	C2	27.841
177	C1	25.555
178	C1	24.264
179	C1	25.516
180	C2v	24.454
181	C2	26.164
182	D3	23.7
183	C1	23.35
184	C1	23.335
185	C1	23.605
186	Cs	24.692
187	C1	24.163
188	C1	23.748
189	C1	23.278
190	C1	23.899
191	C1	24.972
192	C1	24.601
193	C1	25.672
194	C1	24.796
195	C1	24.342
196	C1	23.803
197	C1	24.309
198	C1	24.699
199	C1	24.671
200	C1	23.463
201	C1	23.264
"@

$records = $data -split "`n" | Where-Object { $_.Trim() -ne "" }
$arr = New-Object "object[,]" $records.Count, 2

for ($i = 0; $i -lt $records.Count; $i++) {
    $fields = $records[$i] -split "`t"
    $arr[$i, 0] = $fields[1]
    $arr[$i, 1] = [double]$fields[2]
}

$ws.Range("B10:C201").Value = $arr

# Mirror the author's final selection/scroll state on the C50 tab.
$ws.Activate()
$ws.Range("B202").Select() | Out-Null
